# Decrease every value in column E (rows 2-99), except row 36 which is left
# unchanged (it was not touched by the original edit), by exactly 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)   # column E is the 5th column
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value = $current - 1
    }
}
